$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Theoretical Friction Factor" column (V) values for rows 2-5
# to reflect the difference between predicted and measured temperatures.
$ws.Range("V2").Value = 0.003078677535596517
$ws.Range("V3").Value = 0.002658857872402826
$ws.Range("V4").Value = 0.002048211086161262
$ws.Range("V5").Value = 0.002455551943221542
